$d = $word.ActiveDocument

$replacements = @(
    @("769×5=3845", "274×7=1918"),
    @("691×9=6219", "145×2=290"),
    @("351×9=3159", "654×5=3270"),
    @("241×7=1687", "555×4=2220"),
    @("356×7=2492", "531×7=3717"),
    @("528×2=1056", "621×7=4347"),
    @("269×5=1345", "791×4=3164"),
    @("275×7=1925", "101×2=202"),
    @("778×3=2334", "415×5=2075"),
    @("117×2=234",  "640×6=3840"),
    @("356×8=2848", "302×2=604"),
    @("872×3=2616", "493×4=1972"),
    @("768×3=2304", "447×9=4023"),
    @("491×9=4419", "112×7=784"),
    @("541×7=3787", "627×8=5016"),
    @("902×3=2706", "902×6=5412"),
    @("599×8=4792", "367×8=2936"),
    @("834×3=2502", "190×7=1330"),
    @("326×6=1956", "432×6=2592"),
    @("672×7=4704", "143×8=1144"),
    @("456×6=2736", "565×4=2260"),
    @("668×6=4008", "476×4=1904"),
    @("813×4=3252", "556×4=2224"),
    @("445×7=3115", "231×9=2079"),
    @("740×4=2960", "137×3=411")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
